{"js": "// Replace each \"divided by\" expression in the worksheet table with its\n// new value, per the commit's regenerated numbers. Every original string\n// is unique in the document, so a literal (non-wildcard) body.search()\n// for each old value yields exactly one hit to replace.\nconst replacements = [\n  [\"703\u00f77=\", \"169\u00f77=\"],\n  [\"866\u00f72=\", \"455\u00f76=\"],\n  [\"789\u00f73=\", \"628\u00f73=\"],\n  [\"677\u00f77=\", \"387\u00f73=\"],\n  [\"510\u00f72=\", \"879\u00f76=\"],\n  [\"524\u00f77=\", \"679\u00f76=\"],\n  [\"861\u00f78=\", \"892\u00f78=\"],\n  [\"837\u00f79=\", \"518\u00f72=\"],\n  [\"338\u00f79=\", \"635\u00f76=\"],\n  [\"818\u00f73=\", \"793\u00f72=\"],\n  [\"647\u00f73=\", \"511\u00f72=\"],\n  [\"257\u00f75=\", \"379\u00f75=\"],\n  [\"986\u00f72=\", \"305\u00f75=\"],\n  [\"492\u00f72=\", \"182\u00f72=\"],\n  [\"461\u00f77=\", \"213\u00f74=\"],\n  [\"772\u00f79=\", \"184\u00f79=\"],\n  [\"382\u00f74=\", \"647\u00f78=\"],\n  [\"323\u00f74=\", \"483\u00f74=\"],\n  [\"494\u00f77=\", \"790\u00f75=\"],\n  [\"647\u00f77=\", \"680\u00f78=\"],\n  [\"349\u00f77=\", \"976\u00f77=\"],\n  [\"186\u00f75=\", \"932\u00f72=\"],\n  [\"823\u00f74=\", \"542\u00f79=\"],\n  [\"477\u00f73=\", \"504\u00f76=\"],\n  [\"258\u00f79=\", \"435\u00f76=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const hits = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  hits.load(\"items\");\n  await context.sync();\n\n  if (hits.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (const hit of hits.items) {\n    hit.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"divided by\" expression in the worksheet table with its\n# new value, per the commit's regenerated numbers. Every original string\n# is unique in the document, so Find/Replace (match whole text, no\n# wildcards) on each old value touches exactly one run.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"703\u00f77=\", \"169\u00f77=\"),\n    @(\"866\u00f72=\", \"455\u00f76=\"),\n    @(\"789\u00f73=\", \"628\u00f73=\"),\n    @(\"677\u00f77=\", \"387\u00f73=\"),\n    @(\"510\u00f72=\", \"879\u00f76=\"),\n    @(\"524\u00f77=\", \"679\u00f76=\"),\n    @(\"861\u00f78=\", \"892\u00f78=\"),\n    @(\"837\u00f79=\", \"518\u00f72=\"),\n    @(\"338\u00f79=\", \"635\u00f76=\"),\n    @(\"818\u00f73=\", \"793\u00f72=\"),\n    @(\"647\u00f73=\", \"511\u00f72=\"),\n    @(\"257\u00f75=\", \"379\u00f75=\"),\n    @(\"986\u00f72=\", \"305\u00f75=\"),\n    @(\"492\u00f72=\", \"182\u00f72=\"),\n    @(\"461\u00f77=\", \"213\u00f74=\"),\n    @(\"772\u00f79=\", \"184\u00f79=\"),\n    @(\"382\u00f74=\", \"647\u00f78=\"),\n    @(\"323\u00f74=\", \"483\u00f74=\"),\n    @(\"494\u00f77=\", \"790\u00f75=\"),\n    @(\"647\u00f77=\", \"680\u00f78=\"),\n    @(\"349\u00f77=\", \"976\u00f77=\"),\n    @(\"186\u00f75=\", \"932\u00f72=\"),\n    @(\"823\u00f74=\", \"542\u00f79=\"),\n    @(\"477\u00f73=\", \"504\u00f76=\"),\n    @(\"258\u00f79=\", \"435\u00f76=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
